$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Actualiza base de datos EC": the arrears value ("Valor Mora") that was
# recorded against period 2006 (row 16) actually belongs to period 1909
# (row 25) -- swap the two amounts.
$ws.Range("F16").Value = 33125
$ws.Range("F25").Value = 26500
